$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'317.39"
$ws.Range("E2").Value = "'-2.96%"
$ws.Range("D3").Value = "'41.88"
$ws.Range("E3").Value = "'-5.73%"
$ws.Range("E4").Value = "'1.21%"
$ws.Range("D5").Value = "'0.08071"
$ws.Range("E5").Value = "'-3.69%"
$ws.Range("D6").Value = "'4.378"
$ws.Range("E6").Value = "'-1.62%"
$ws.Range("D7").Value = "'1.753"
$ws.Range("E7").Value = "'-9.52%"
$ws.Range("D8").Value = "'0.9280"
$ws.Range("E8").Value = "'-4.49%"
$ws.Range("D9").Value = "'0.1124"
$ws.Range("E9").Value = "'-0.24%"
$ws.Range("D10").Value = "'0.1855"
$ws.Range("E10").Value = "'-2.39%"
$ws.Range("D11").Value = "'0.09256"
$ws.Range("E11").Value = "'-4.26%"
$ws.Range("D12").Value = "'0.04563"
$ws.Range("E12").Value = "'-1.47%"
$ws.Range("D13").Value = "'7.371"
$ws.Range("E13").Value = "'-15.57%"
$ws.Range("D15").Value = "'0.001293"
$ws.Range("E15").Value = "'0.59%"
$ws.Range("D16").Value = "'0.005973"
$ws.Range("E16").Value = "'3.10%"
$ws.Range("D17").Value = "'3.355"
$ws.Range("E17").Value = "'-1.45%"
$ws.Range("D19").Value = "'0.3398"
$ws.Range("E19").Value = "'1.09%"
$ws.Range("E20").Value = "'1.31%"
$ws.Range("D21").Value = "'0.2605"
$ws.Range("E21").Value = "'0.81%"
$ws.Range("D22").Value = "'0.04174"
$ws.Range("E22").Value = "'0.11%"
$ws.Range("D23").Value = "'0.001243"
$ws.Range("E23").Value = "'0.07%"
$ws.Range("D24").Value = "'0.004316"
$ws.Range("E24").Value = "'-2.35%"
$ws.Range("D25").Value = "'0.0001222"
$ws.Range("E25").Value = "'-6.36%"
$ws.Range("E26").Value = "'-0.09%"
$ws.Range("D38").Value = "'0.02574"
$ws.Range("E38").Value = "'-5.75%"
$ws.Range("D39").Value = "'0.05440"
$ws.Range("E39").Value = "'-3.07%"
$ws.Range("D40").Value = "'0.008075"
$ws.Range("E40").Value = "'2.63%"
$ws.Range("D41").Value = "'0.1391"
$ws.Range("E41").Value = "'-1.21%"
$ws.Range("D42").Value = "'0.007556"
$ws.Range("E42").Value = "'2.42%"
$ws.Range("D43").Value = "'0.002084"
$ws.Range("E43").Value = "'-1.55%"
$ws.Range("D44").Value = "'0.008241"
$ws.Range("E44").Value = "'3.97%"
$ws.Range("D45").Value = "'0.3147"
$ws.Range("E45").Value = "'-9.92%"
$ws.Range("D46").Value = "'0.00006801"
$ws.Range("E46").Value = "'-1.92%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("D48").Value = "'0.003391"
$ws.Range("E48").Value = "'-3.20%"
$ws.Range("D49").Value = "'0.004112"
$ws.Range("E49").Value = "'16.05%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'-0.20%"
